# Add 5 new data rows (14-18) to Sheet1, matching the uploaded source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 1010
$ws.Range("B14").Value = "test 10"
$ws.Range("C14").Value = "j"
$ws.Range("D14").Value = "PRJ-09"

# Row 15
$ws.Range("A15").Value = 1003
$ws.Range("B15").Value = "Test Dup"
$ws.Range("C15").Value = 152
$ws.Range("D15").Value = "PRJ-02"

# Row 16
$ws.Range("A16").Value = 1011
$ws.Range("B16").Value = "Test Coorect"
$ws.Range("C16").Value = 1520
$ws.Range("D16").Value = "PRJ-01"

# Row 17
$ws.Range("A17").Value = 1012
$ws.Range("B17").Value = "gf"
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = "PRJ-04"

# Row 18
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "test"
$ws.Range("C18").Value = 155
$ws.Range("D18").Value = "PRJ-01"

# Match the final selection left behind in the source workbook.
$ws.Range("D18").Select()
